$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.698.24'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '3.075.54'
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.071.28'
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("E9").Value = '  -2.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("E14").Value = '  -3.81%  '
$ws.Range("E15").Value = '  -1.77%  '
$ws.Range("D16").Value = '3.584.73'
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").Value = '66.703.16'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.48%  '
$ws.Range("E19").Value = '  -2.49%  '
$ws.Range("D20").Value = '3.076.55'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '488.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.49%  '
$ws.Range("E22").Value = '  -2.68%  '
$ws.Range("E23").Value = '  -3.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.18%  '
$ws.Range("E26").Value = '  -3.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.82'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("E30").Value = '  -4.57%  '
$ws.Range("E31").Value = '  -2.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.43%  '
$ws.Range("E33").Value = '  -3.39%  '
$ws.Range("D34").Value = '0.0₃0913'
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.946'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.87%  '
$ws.Range("B38").Value = 'Arweave'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '46.95'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.301'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.89%  '
$ws.Range("D43").Value = '2.761.56'
$ws.Range("E43").Value = '  -2.39%  '
$ws.Range("E44").Value = '  -2.40%  '
$ws.Range("E45").Value = '  -3.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '369.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '135.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.69'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("E51").Value = '  -2.09%  '
